# Insert a new data row above row 61 (pushing existing rows 61..178 down to 62..179)
# and populate it with a new weekly price record for "Betarraga" at
# "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("61:61").Insert()

$ws.Range("A61").Value = 5
$ws.Range("B61").Value = "Macroferia Regional de Talca"
$ws.Range("C61").Value = "Maule"
$ws.Range("D61").Value = 44469
$ws.Range("E61").Value = 7
$ws.Range("F61").Value = 100114014
$ws.Range("G61").Value = "Betarraga"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 650
$ws.Range("L61").Value = 650
$ws.Range("M61").Value = 650
$ws.Range("N61").Value = "$/paquete 5 unidades"
$ws.Range("O61").Value = "Región del Maule"
$ws.Range("P61").Value = 130
$ws.Range("Q61").Value = 5
$ws.Range("R61").Value = "Hortaliza"
